$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 2 (old rows 2-6 shift down to 5-9)
$ws.Range("A2:A4").EntireRow.Insert()

# Row 2: BEFORE_REPLACEMENT, B2 left empty
$ws.Range("A2").Value = "BEFORE_REPLACEMENT"

# Row 3: AFTER_REPLACEMENT, B3 = single space
$ws.Range("A3").Value = "AFTER_REPLACEMENT"
$ws.Range("B3").Value = " "

# Row 4: LINK_EDIT_FILE
$ws.Range("A4").Value = "LINK_EDIT_FILE"
$ws.Range("B4").Value = "https://docs.google.com/spreadsheets/d/16uVFfVMKR7jVXA70g4BCo8KAE7iZVYnJT48oTpD1Z-4/edit?gid=0#gid=00"

# Update old "ddal" row (now row 7) B cell to add a trailing newline
$ws.Range("B7").Value = "- Uống thuốc theo toa. Tái khám theo hẹn hoặc khi có dấu hiệu bất thường.`n- Ăn lạt, nhiều rau xanh, ngũ cốc nguyên hạt…`n- Hạn chế ăn thịt đỏ, nội tạng động vật, mỡ, trứng,tinh bột`n"
$ws.Rows.Item(7).AutoFit()

# New row 10 (B10 stored as text so it reads back as "123312", not a number)
$ws.Range("A10").Value = "testso"
$ws.Range("B10").Value = "'123312"

$ws.Calculate()
